$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Major Update to OOP Structure, Pending Requests": add the new
# num_assigned (column D) counters for every faculty row (2-76),
# initialised to 0 pending requests.
$ws.Range("D2:D76").Value = 0

# Leave the sheet with that newly-populated range selected (mirrors the
# author's final on-screen selection in the bottom-left frozen pane).
$ws.Range("D2:D76").Select() | Out-Null
